$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.280.03'
$ws.Range("D3").Value = '1.667.70'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07825'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '1.639.87'
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").Value = '1.894.57'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5533'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.72'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.681'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.038'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.010'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1227'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.198'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.485'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05887'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.279'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.604'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.278'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.613'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9621'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.422'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5802'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01608'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8637'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.848'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.049.70'
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.009'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").Value = '1.804.96'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -4.16%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.012'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4378'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.008'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05160'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.429'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.11%  '
